$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove rows 5 and 6 ---
# This shifts the original row 7 (which already carries the "thick bottom border" row
# formatting) up into row 5's position, so the resulting last row keeps that bottom-border
# row styling. Rows 1-4 are left completely untouched by this operation.
$ws.Range("A5:G6").EntireRow.Delete()

# --- Normalize the D:F styling on the (new) row 5 to match the rest of the row (style of A2:C2) ---
$ws.Range("A2:C2").Copy()
$ws.Range("D5").PasteSpecial(-4122)  # xlPasteFormats

# --- Update header titles (row 1) ---
$ws.Range("A1").Value = "16.4.2.1 Алынган жана өз эрки менен берилген ок атуучу куралдардын саны "
$ws.Range("C1").Value = "16.4.2.1 Number of seized and voluntary surrendered firearms"

# --- Add the new year column (H) ---
$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H3").Value = 2020

$ws.Range("G4").Copy()
$ws.Range("H4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H4").Value = 158

$ws.Range("G5").Copy()
$ws.Range("H5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H5").Value = 397

# --- Update row 4 text/values ---
$ws.Range("A4").Value = "Өз ыктыяры менен тапшырган ок атуучу куралдар"
$ws.Range("C4").Value = "Voluntarily surrendered firearms"
$ws.Range("G4").Value = 146

# --- Update row 5 text/values ---
$ws.Range("A5").Value = "Алынган ок атуучу куралдар"
$ws.Range("B5").Value = "Изъятые огнестрельные оружия"
$ws.Range("C5").Value = "Seized firearms"
$ws.Range("D5").Value = 217
$ws.Range("E5").Value = 399
$ws.Range("F5").Value = 296
$ws.Range("G5").Value = 127

# --- Reset selection to A1 (closest achievable state to a cleared selection) ---
$ws.Range("A1").Select()
